# Update NATMI LR-pair output values with refreshed TPM-derived numbers.
# The underlying row/column layout is unchanged; only the computed
# statistics in columns E:T for data rows 2-10 move to new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.007630333333333334
$ws.Range("H2").Value = 0.022891
$ws.Range("I2").Value = 0.005544871895800688
$ws.Range("J2").Value = 0.005544871895800688
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003710666666666666
$ws.Range("N2").Value = 0.011132
$ws.Range("O2").Value = 0.001642024256586498
$ws.Range("P2").Value = 0.001642024256586498
$ws.Range("Q2").Value = 0.00002831362355555555
$ws.Range("R2").Value = 0.000254822612
$ws.Range("S2").Value = 0.000009104814152569492
$ws.Range("T2").Value = 0.000009104814152569492

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.007630333333333334
$ws.Range("H3").Value = 0.022891
$ws.Range("I3").Value = 0.005544871895800688
$ws.Range("J3").Value = 0.005544871895800688
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.049608666666666
$ws.Range("N3").Value = 6.148826
$ws.Range("O3").Value = 0.9069818039462568
$ws.Range("P3").Value = 0.9069818039462569
$ws.Range("Q3").Value = 0.01563919732955555
$ws.Range("R3").Value = 0.140752775966
$ws.Range("S3").Value = 0.005029097914704209
$ws.Range("T3").Value = 0.00502909791470421

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.007630333333333334
$ws.Range("H4").Value = 0.022891
$ws.Range("I4").Value = 0.005544871895800688
$ws.Range("J4").Value = 0.005544871895800688
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.206493
$ws.Range("N4").Value = 0.619479
$ws.Range("O4").Value = 0.09137617179715662
$ws.Range("P4").Value = 0.09137617179715661
$ws.Range("Q4").Value = 0.001575610421
$ws.Range("R4").Value = 0.014180493789
$ws.Range("S4").Value = 0.0005066691669439092
$ws.Range("T4").Value = 0.0005066691669439091

$ws.Range("G5").Value = 0.4678513333333334
$ws.Range("I5").Value = 0.3399819636031033
$ws.Range("J5").Value = 0.3399819636031033
$ws.Range("M5").Value = 0.003710666666666666
$ws.Range("N5").Value = 0.011132
$ws.Range("O5").Value = 0.001642024256586498
$ws.Range("P5").Value = 0.001642024256586498
$ws.Range("Q5").Value = 0.001736040347555556
$ws.Range("R5").Value = 0.015624363128
$ws.Range("S5").Value = 0.0005582586310382037
$ws.Range("T5").Value = 0.0005582586310382037

$ws.Range("G6").Value = 0.4678513333333334
$ws.Range("I6").Value = 0.3399819636031033
$ws.Range("J6").Value = 0.3399819636031033
$ws.Range("M6").Value = 2.049608666666666
$ws.Range("N6").Value = 6.148826
$ws.Range("O6").Value = 0.9069818039462568
$ws.Range("P6").Value = 0.9069818039462569
$ws.Range("Q6").Value = 0.9589121475115555
$ws.Range("R6").Value = 8.630209327604002
$ws.Range("S6").Value = 0.3083574546579333
$ws.Range("T6").Value = 0.3083574546579333

$ws.Range("G7").Value = 0.4678513333333334
$ws.Range("I7").Value = 0.3399819636031033
$ws.Range("J7").Value = 0.3399819636031033
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.206493
$ws.Range("N7").Value = 0.619479
$ws.Range("O7").Value = 0.09137617179715662
$ws.Range("P7").Value = 0.09137617179715661
$ws.Range("Q7").Value = 0.09660802537400001
$ws.Range("R7").Value = 0.8694722283660001
$ws.Range("S7").Value = 0.03106625031413182
$ws.Range("T7").Value = 0.03106625031413182

$ws.Range("G8").Value = 0.9006246666666667
$ws.Range("H8").Value = 2.701874
$ws.Range("I8").Value = 0.6544731645010959
$ws.Range("J8").Value = 0.6544731645010959
$ws.Range("M8").Value = 0.003710666666666666
$ws.Range("N8").Value = 0.011132
$ws.Range("O8").Value = 0.001642024256586498
$ws.Range("P8").Value = 0.001642024256586498
$ws.Range("Q8").Value = 0.003341917929777778
$ws.Range("R8").Value = 0.030077261368
$ws.Range("S8").Value = 0.001074660811395725
$ws.Range("T8").Value = 0.001074660811395725

$ws.Range("G9").Value = 0.9006246666666667
$ws.Range("H9").Value = 2.701874
$ws.Range("I9").Value = 0.6544731645010959
$ws.Range("J9").Value = 0.6544731645010959
$ws.Range("M9").Value = 2.049608666666666
$ws.Range("N9").Value = 6.148826
$ws.Range("O9").Value = 0.9069818039462568
$ws.Range("P9").Value = 0.9069818039462569
$ws.Range("Q9").Value = 1.845928122213778
$ws.Range("R9").Value = 16.613353099924
$ws.Range("S9").Value = 0.5935952513736192
$ws.Range("T9").Value = 0.5935952513736193

$ws.Range("G10").Value = 0.9006246666666667
$ws.Range("H10").Value = 2.701874
$ws.Range("I10").Value = 0.6544731645010959
$ws.Range("J10").Value = 0.6544731645010959
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.206493
$ws.Range("N10").Value = 0.619479
$ws.Range("O10").Value = 0.09137617179715662
$ws.Range("P10").Value = 0.09137617179715661
$ws.Range("Q10").Value = 0.185972689294
$ws.Range("R10").Value = 1.673754203646
$ws.Range("S10").Value = 0.05980325231608089
$ws.Range("T10").Value = 0.05980325231608088

